$d = $word.ActiveDocument

# Step 1: narrow range for the part before where "005" should go, replace with itself (no-op) just to test
$rA = $d.Range(37, 52)   # "Volunteer # ___"
Write-Output "A=[$($rA.Text)]"
$rA.Find.Execute("Volunteer # ___", $false,$false,$false,$false,$false,$true,1,$false,"Volunteer # ___",2)

$rB = $d.Range(52, 52)
$rB.InsertAfter("005")
